$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Modules_and_forms")

# Insert two new columns (E and F) for case labels
$ws.Columns("E:F").Insert()

# Set column widths for the newly inserted columns
# (the stored worksheet XML width is ColumnWidth + 5/6, so compensate here
# so the saved width matches 19.5 / 17.5)
$ws.Columns("E").ColumnWidth = 19.5 - (5/6)
$ws.Columns("F").ColumnWidth = 17.5 - (5/6)

# Populate the header row
$ws.Range("E1").Value = "label_for_cases_en"
$ws.Range("F1").Value = "label_for_cases_fra"

# Populate the "Cases" label for the Module row
$ws.Range("E2").Value = "Cases"
$ws.Range("F2").Value = "Cases"

# Make this sheet the active sheet/tab and set its selection
$ws.Activate() | Out-Null
$ws.Range("F4").Select() | Out-Null
